# New .csv- and .xlsx-data files with subsample n = 41
# Trims the data set from 51 subjects (rows 2-52) down to 41 subjects
# (rows 2-42), re-labels the WPM_log/FPM_log columns as FPM/WPM, and
# replaces the log-transformed D/E values with the newly resampled ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the trailing 10 subjects (rows 43-52) -----------------------
$ws.Rows("43:52").Delete() | Out-Null

# --- Re-label the D/E headers ------------------------------------------
$ws.Range("D1").Value = "FPM"
$ws.Range("E1").Value = "WPM"

# --- New column D ("FPM") values for rows 2-42 --------------------------
$dVals = @(
    6.952134847798976,
    7.491523672377968,
    6.949746855268403,
    6.983665645428913,
    7.445685204605836,
    7.139870140659839,
    7.37077782435057,
    7.308278077087381,
    6.842542813453106,
    6.974849819532158,
    7.078653115998773,
    6.634880798133581,
    7.108065594295618,
    6.730530529381854,
    7.171153351207309,
    7.146766028094526,
    7.090863241949512,
    7.291442566904449,
    7.314461508176811,
    7.132677768078501,
    7.189083292612422,
    7.233796984705399,
    7.104731651556778,
    6.980222777484832,
    6.964743475718091,
    7.350428595077923,
    6.858542388421064,
    7.174472723554255,
    6.913995806023228,
    7.057897253003959,
    6.770001682661413,
    6.514652591241854,
    7.110491622815228,
    7.137692870403487,
    7.058374451406982,
    7.333824634405921,
    7.105579474676065,
    7.476156436258767,
    7.010305555679961,
    7.02437849804367,
    7.18041331328177
)

# --- New column E ("WPM") values for rows 2-42 ---------------------------
$eVals = @(
    6.783413977560148,
    7.270794497089682,
    6.704910435309077,
    6.635059784800578,
    7.24539013591141,
    6.920435962115727,
    7.089231582809733,
    7.071794403651529,
    6.678030764716598,
    6.920491249657466,
    6.869212198532422,
    6.364662321020626,
    6.937830853335288,
    6.464115853763414,
    6.918771391529022,
    6.867175612586004,
    6.93695694447592,
    7.180032388973652,
    6.998034422000302,
    6.845700061654469,
    6.962926249785784,
    7.023767386392288,
    6.822219003719034,
    6.615899390546929,
    6.728627430655716,
    7.126630978545323,
    6.61421054518505,
    7.095960151653968,
    6.749686950380546,
    6.774180426247743,
    6.54429203562524,
    6.404263772480118,
    6.93005144194609,
    6.902993022545933,
    6.927850423299475,
    7.176516879949531,
    6.821838095843244,
    7.303612586750316,
    6.821495995747958,
    6.726143734894434,
    7.035865229496709
)

for ($i = 0; $i -lt $dVals.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}
